$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.807.68'
$ws.Range('E2').Value = '  +7.05%  '

$ws.Range('D3').Value = '3.308.66'
$ws.Range('E3').Value = '  +2.49%  '

$ws.Range('E4').Value = '  +0.33%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '409.35'
$ws.Range('E5').Value = '  +3.99%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '111.93'
$ws.Range('E6').Value = '  +4.48%  '

$ws.Range('D7').Value = '3.301.08'
$ws.Range('E7').Value = '  +2.33%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.564'
$ws.Range('E8').Value = '  -1.61%  '

$ws.Range('E9').Value = '  +0.08%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.615'
$ws.Range('E10').Value = '  +0.53%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.107'
$ws.Range('E11').Value = '  +11.74%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '38.59'
$ws.Range('E12').Value = '  -1.14%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.142'
$ws.Range('E13').Value = '  +0.51%  '

$ws.Range('D14').Value = '3.858.70'
$ws.Range('E14').Value = '  +3.21%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.15'
$ws.Range('E15').Value = '  +0.32%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.00'
$ws.Range('E16').Value = '  -0.01%  '

$ws.Range('D17').Value = '3.389.23'
$ws.Range('E17').Value = '  +5.08%  '

$ws.Range('D18').Value = '60.861.20'
$ws.Range('E18').Value = '  +7.46%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.991'
$ws.Range('E19').Value = '  -3.91%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.39'
$ws.Range('E20').Value = '  -5.17%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000111'
$ws.Range('E21').Value = '  +5.52%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.22'
$ws.Range('E22').Value = '  -3.45%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '293.45'
$ws.Range('E23').Value = '  -1.22%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.16'
$ws.Range('E24').Value = '  -6.47%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.16'
$ws.Range('E25').Value = '  -1.05%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.06'
$ws.Range('E26').Value = '  -2.44%  '

$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.47'
$ws.Range('E27').Value = '  +2.73%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '28.40'
$ws.Range('E28').Value = '  +2.11%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.39'
$ws.Range('E29').Value = '  +2.42%  '

$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.168'
$ws.Range('E30').Value = '  -0.64%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.41'
$ws.Range('E31').Value = '  -3.39%  '

$ws.Range('B32').Value = 'Dai'
$ws.Range('C32').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.998'
$ws.Range('E32').Value = '  -0.16%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.109'
$ws.Range('E33').Value = '  +0.31%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.13'
$ws.Range('E34').Value = '  -1.51%  '

$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '40.10'
$ws.Range('E35').Value = '  +6.62%  '

$ws.Range('B36').Value = 'Toncoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.39'
$ws.Range('E36').Value = '  +13.16%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0471'
$ws.Range('E37').Value = '  -2.51%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '51.93'
$ws.Range('E38').Value = '  +0.45%  '

$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.24%  '

$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.08'
$ws.Range('E40').Value = '  +3.16%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.29'
$ws.Range('E41').Value = '  -6.99%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '137.01'
$ws.Range('E42').Value = '  +2.11%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.118'
$ws.Range('E43').Value = '  -1.50%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.86'
$ws.Range('E44').Value = '  -1.32%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.275'
$ws.Range('E45').Value = '  -2.19%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.76'
$ws.Range('E46').Value = '  -5.60%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '16.14'
$ws.Range('E47').Value = '  -4.96%  '

$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.16'
$ws.Range('E48').Value = '  +2.56%  '

$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '3.725.42'
$ws.Range('E49').Value = '  +4.99%  '

$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.126.14'
$ws.Range('E50').Value = '  -1.09%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.73'
$ws.Range('E51').Value = '  -6.15%  '

Write-Host "Updated cryptos list"